$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -19.19123757349545
$ws.Range("C2").Value = 2.521691129046427
$ws.Range("D2").Value = -19.19123757349545
$ws.Range("E2").Value = -19.19123757349545
$ws.Range("F2").Value = -19.19123757349545
$ws.Range("G2").Value = -19.19123757349545
$ws.Range("H2").Value = -19.19123757349545
$ws.Range("I2").Value = -19.19123757349545
$ws.Range("J2").Value = -19.19123757349545
$ws.Range("K2").Value = -19.19123757349545
$ws.Range("B3").Value = -19.19123757349545
$ws.Range("C3").Value = -19.19123757349545
$ws.Range("D3").Value = -19.19123757349545
$ws.Range("E3").Value = -19.19123757349545
$ws.Range("F3").Value = -19.19123757349545
$ws.Range("G3").Value = -19.19123757349545
$ws.Range("H3").Value = -19.19123757349545
$ws.Range("I3").Value = -19.19123757349545
$ws.Range("J3").Value = -19.19123757349545
$ws.Range("K3").Value = -19.19123757349545
$ws.Range("B4").Value = -19.19123757349545
$ws.Range("C4").Value = 2.20646264677666
$ws.Range("D4").Value = 2.928118058899412
$ws.Range("E4").Value = -19.19123757349545
$ws.Range("F4").Value = 2.539518209220593
$ws.Range("G4").Value = -19.19123757349545
$ws.Range("H4").Value = 1.873908330344156
$ws.Range("I4").Value = -19.19123757349545
$ws.Range("J4").Value = 2.43217823179939
$ws.Range("K4").Value = -19.19123757349545
$ws.Range("B5").Value = -19.19123757349545
$ws.Range("C5").Value = 0.9665694574430946
$ws.Range("D5").Value = -19.19123757349545
$ws.Range("E5").Value = -19.19123757349545
$ws.Range("F5").Value = -19.19123757349545
$ws.Range("G5").Value = 2.217583999872494
$ws.Range("H5").Value = -19.19123757349545
$ws.Range("I5").Value = -19.19123757349545
$ws.Range("J5").Value = -19.19123757349545
$ws.Range("K5").Value = -19.19123757349545
$ws.Range("B6").Value = -19.19123757349545
$ws.Range("C6").Value = -19.19123757349545
$ws.Range("D6").Value = -19.19123757349545
$ws.Range("E6").Value = -19.19123757349545
$ws.Range("F6").Value = -19.19123757349545
$ws.Range("G6").Value = -19.19123757349545
$ws.Range("H6").Value = -19.19123757349545
$ws.Range("I6").Value = -19.19123757349545
$ws.Range("J6").Value = -19.19123757349545
$ws.Range("K6").Value = -19.19123757349545
$ws.Range("B7").Value = 2.979790446868304
$ws.Range("C7").Value = -19.19123757349545
$ws.Range("D7").Value = -19.19123757349545
$ws.Range("E7").Value = -19.19123757349545
$ws.Range("F7").Value = -19.19123757349545
$ws.Range("G7").Value = -19.19123757349545
$ws.Range("H7").Value = -19.19123757349545
$ws.Range("I7").Value = -19.19123757349545
$ws.Range("J7").Value = -19.19123757349545
$ws.Range("K7").Value = -19.19123757349545
$ws.Range("B8").Value = -19.19123757349545
$ws.Range("C8").Value = -19.19123757349545
$ws.Range("D8").Value = -19.19123757349545
$ws.Range("E8").Value = 2.905870277490526
$ws.Range("F8").Value = -19.19123757349545
$ws.Range("G8").Value = -19.19123757349545
$ws.Range("H8").Value = -19.19123757349545
$ws.Range("I8").Value = -19.19123757349545
$ws.Range("J8").Value = -19.19123757349545
$ws.Range("K8").Value = -19.19123757349545
$ws.Range("B9").Value = 3.598276370873494
$ws.Range("C9").Value = -19.19123757349545
$ws.Range("D9").Value = -19.19123757349545
$ws.Range("E9").Value = -19.19123757349545
$ws.Range("F9").Value = -19.19123757349545
$ws.Range("G9").Value = -19.19123757349545
$ws.Range("H9").Value = -19.19123757349545
$ws.Range("I9").Value = -19.19123757349545
$ws.Range("J9").Value = -19.19123757349545
$ws.Range("K9").Value = -19.19123757349545
$ws.Range("B10").Value = -19.19123757349545
$ws.Range("C10").Value = -19.19123757349545
$ws.Range("D10").Value = -19.19123757349545
$ws.Range("E10").Value = -19.19123757349545
$ws.Range("F10").Value = -19.19123757349545
$ws.Range("G10").Value = -19.19123757349545
$ws.Range("H10").Value = -19.19123757349545
$ws.Range("I10").Value = -19.19123757349545
$ws.Range("J10").Value = -19.19123757349545
$ws.Range("K10").Value = 2.197642131326067
$ws.Range("B11").Value = -19.19123757349545
$ws.Range("C11").Value = -19.19123757349545
$ws.Range("D11").Value = -19.19123757349545
$ws.Range("E11").Value = 1.970510025574399
$ws.Range("F11").Value = -19.19123757349545
$ws.Range("G11").Value = 2.665989062730902
$ws.Range("H11").Value = -19.19123757349545
$ws.Range("I11").Value = -19.19123757349545
$ws.Range("J11").Value = -19.19123757349545
$ws.Range("K11").Value = 1.402081507582911
$ws.Range("B12").Value = -19.19123757349545
$ws.Range("C12").Value = -19.19123757349545
$ws.Range("D12").Value = -19.19123757349545
$ws.Range("E12").Value = -19.19123757349545
$ws.Range("F12").Value = -19.19123757349545
$ws.Range("G12").Value = -19.19123757349545
$ws.Range("H12").Value = -19.19123757349545
$ws.Range("I12").Value = -19.19123757349545
$ws.Range("J12").Value = -19.19123757349545
$ws.Range("K12").Value = -19.19123757349545
$ws.Range("B13").Value = -19.19123757349545
$ws.Range("C13").Value = -19.19123757349545
$ws.Range("D13").Value = -19.19123757349545
$ws.Range("E13").Value = 1.64994359961846
$ws.Range("F13").Value = -19.19123757349545
$ws.Range("G13").Value = -19.19123757349545
$ws.Range("H13").Value = -19.19123757349545
$ws.Range("I13").Value = -19.19123757349545
$ws.Range("J13").Value = 2.203095845700719
$ws.Range("K13").Value = 1.709874588968424
$ws.Range("B14").Value = -19.19123757349545
$ws.Range("C14").Value = -19.19123757349545
$ws.Range("D14").Value = 1.666346359873598
$ws.Range("E14").Value = -19.19123757349545
$ws.Range("F14").Value = -19.19123757349545
$ws.Range("G14").Value = -19.19123757349545
$ws.Range("H14").Value = -19.19123757349545
$ws.Range("I14").Value = -19.19123757349545
$ws.Range("J14").Value = -19.19123757349545
$ws.Range("K14").Value = 2.086397747849883
$ws.Range("B15").Value = -19.19123757349545
$ws.Range("C15").Value = -19.19123757349545
$ws.Range("D15").Value = -0.2239859803426696
$ws.Range("E15").Value = -19.19123757349545
$ws.Range("F15").Value = -19.19123757349545
$ws.Range("G15").Value = -19.19123757349545
$ws.Range("H15").Value = -19.19123757349545
$ws.Range("I15").Value = -19.19123757349545
$ws.Range("J15").Value = -19.19123757349545
$ws.Range("K15").Value = -19.19123757349545
$ws.Range("B16").Value = -19.19123757349545
$ws.Range("C16").Value = -19.19123757349545
$ws.Range("D16").Value = -19.19123757349545
$ws.Range("E16").Value = -19.19123757349545
$ws.Range("F16").Value = -19.19123757349545
$ws.Range("G16").Value = -19.19123757349545
$ws.Range("H16").Value = -19.19123757349545
$ws.Range("I16").Value = -19.19123757349545
$ws.Range("J16").Value = 2.269692705033327
$ws.Range("K16").Value = -19.19123757349545
$ws.Range("B17").Value = -19.19123757349545
$ws.Range("C17").Value = 0.7477575970324786
$ws.Range("D17").Value = -0.07276029003023325
$ws.Range("E17").Value = -19.19123757349545
$ws.Range("F17").Value = -19.19123757349545
$ws.Range("G17").Value = -19.19123757349545
$ws.Range("H17").Value = 0.7440210147734231
$ws.Range("I17").Value = -19.19123757349545
$ws.Range("J17").Value = 1.428879312243232
$ws.Range("K17").Value = -19.19123757349545
$ws.Range("B18").Value = -19.19123757349545
$ws.Range("C18").Value = -19.19123757349545
$ws.Range("D18").Value = -19.19123757349545
$ws.Range("E18").Value = -19.19123757349545
$ws.Range("F18").Value = -19.19123757349545
$ws.Range("G18").Value = -19.19123757349545
$ws.Range("H18").Value = 0.348612257423252
$ws.Range("I18").Value = -19.19123757349545
$ws.Range("J18").Value = 1.312464930588257
$ws.Range("K18").Value = -19.19123757349545
$ws.Range("B19").Value = -19.19123757349545
$ws.Range("C19").Value = -19.19123757349545
$ws.Range("D19").Value = 1.670334784838688
$ws.Range("E19").Value = -19.19123757349545
$ws.Range("F19").Value = -19.19123757349545
$ws.Range("G19").Value = -19.19123757349545
$ws.Range("H19").Value = 2.024928632403775
$ws.Range("I19").Value = 4.321925805283625
$ws.Range("J19").Value = -19.19123757349545
$ws.Range("K19").Value = -19.19123757349545
$ws.Range("B20").Value = -19.19123757349545
$ws.Range("C20").Value = 1.523411264088333
$ws.Range("D20").Value = 2.078887194036384
$ws.Range("E20").Value = -19.19123757349545
$ws.Range("F20").Value = 3.826398160961791
$ws.Range("G20").Value = -19.19123757349545
$ws.Range("H20").Value = 2.12597241646091
$ws.Range("I20").Value = -19.19123757349545
$ws.Range("J20").Value = -19.19123757349545
$ws.Range("K20").Value = 2.392776647271366
$ws.Range("B21").Value = -19.19123757349545
$ws.Range("C21").Value = 1.647976212091085
$ws.Range("D21").Value = -19.19123757349545
$ws.Range("E21").Value = 2.445725510682597
$ws.Range("F21").Value = -19.19123757349545
$ws.Range("G21").Value = 3.170279243718932
$ws.Range("H21").Value = 2.307974274738089
$ws.Range("I21").Value = -19.19123757349545
$ws.Range("J21").Value = -19.19123757349545
$ws.Range("K21").Value = -19.19123757349545
